# The "vocabulary" sheet had an extra, now-unneeded helper column G
# (it only held a couple of stray description/unit notes such as
# "extra description of the factor levels below", "millimol",
# "micromol", "testing purpose"). Remove the whole column, which
# shifts every column from H onward one position to the left.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(7).EntireColumn.Delete() | Out-Null

# The default "list" data validation on the default_unit row referenced
# the (now shifted) helper columns; point it at the new location.
$ws.Range("D6:F6").Validation.Modify(3, 1, 1, "`$G`$6:`$DB`$6") | Out-Null

# Leave the sheet focused on the cell the author ended up looking at.
$ws.Activate() | Out-Null
$ws.Range("F5").Select() | Out-Null
